$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns keep their original text formatting
# when we write the new values, so Excel does not reinterpret them as numbers/percentages.
$updates = @{
    2 = @{ D = "287.46"; E = "1.95%" }
    3 = @{ D = "29.19"; E = "3.29%" }
    4 = @{ D = "5.178"; E = "2.45%" }
    5 = @{ D = "0.06938"; E = "6.98%" }
    6 = @{ D = "7.392"; E = "1.75%" }
    7 = @{ D = "3.550"; E = "5.19%" }
    8 = @{ D = "1.402"; E = "2.15%" }
    9 = @{ D = "0.8987"; E = "-3.30%" }
    10 = @{ D = "0.1597"; E = "2.66%" }
    11 = @{ D = "0.07329"; E = "27.56%" }
    12 = @{ D = "0.07646"; E = "1.12%" }
    13 = @{ D = "0.02930"; E = "1.58%" }
    14 = @{ D = "0.08989"; E = "0.02%" }
    15 = @{ D = "0.001589"; E = "0.15%" }
    16 = @{ D = "0.0006452"; E = "1.38%" }
    17 = @{ D = "0.006346"; E = "4.93%" }
    18 = @{ D = "3.458"; E = "0.22%" }
    19 = @{ D = "2.227"; E = "-0.31%" }
    20 = @{ D = "0.3203"; E = "0.08%" }
    21 = @{ D = "0.1324"; E = "1.63%" }
    22 = @{ D = "4.014"; E = "-2.00%" }
    23 = @{ D = "0.1554"; E = "1.72%" }
    24 = @{ D = "0.04527"; E = "1.28%" }
    25 = @{ D = "0.001211"; E = "2.20%" }
    26 = @{ D = "0.004362"; E = "-0.44%" }
    27 = @{ D = "0.0001170"; E = "-6.35%" }
    28 = @{ D = "0.0001617"; E = "-0.06%" }
    40 = @{ D = "0.04362"; E = "5.05%" }
    41 = @{ D = "0.006935"; E = "4.77%" }
    42 = @{ D = "0.1242"; E = "1.74%" }
    43 = @{ D = "0.002080"; E = "3.02%" }
    44 = @{ D = "0.01182"; E = "-2.19%" }
    45 = @{ D = "0.00005808"; E = "4.83%" }
    47 = @{ D = "0.01307"; E = "0.48%" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $dCell.NumberFormat = "@"
    $eCell.NumberFormat = "@"
    $dCell.Value = $vals.D
    $eCell.Value = $vals.E
}
